$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for Wins/Losses/Ties, matching formatting of existing header row
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy the style of an existing header cell (Z1) onto the new header cells
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("AC1:AE1").PasteSpecial(-4122) | Out-Null

$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 53
    $ws.Cells.Item($r, 30).Value = 61
    $ws.Cells.Item($r, 31).Value = 0
}

Write-Host "done"
